# contactDetail upload sheet - final commit edits
# - corrects a handful of contact fields (names, email, hobbies list, street)
# - tidies up formatting: makes the Pincode/Phone number columns render in solid
#   black (was theme-based "automatic" black) and bumps the header/data row
#   height slightly (18.75pt -> 19.5pt) to match the refreshed layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (first contact) ---------------------------------------------
$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"

# --- Row 3 (second contact) ---------------------------------------------
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Row heights: 18.75pt -> 19.5pt for the header + both data rows ------
$ws.Rows("1:3").RowHeight = 19.5

# --- Pincode / Phone columns: force an explicit black font color ---------
$ws.Range("I2:I3").Font.Color = 0
$ws.Range("K2:K3").Font.Color = 0
